# Updates cryptos list values (price + 1h volume %) per the
# "Updated cryptos list ... with GitHub Actions" commit, including
# the row-49/row-50 swap (Cosmos <-> VeChain).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price strings that look like plain numbers (e.g. "5.50") would be
# auto-converted to numeric values by a direct .Value assignment,
# which truncates meaningful trailing zeros / exact text (e.g. "5.50"
# -> 5.5). To keep them as literal text (matching the source data,
# which stores every Price/Volume cell as a string), such cells are
# written as a literal-string formula and then converted in place to
# a plain value via Copy + PasteSpecial(xlPasteValues), which keeps
# the exact text without touching the cells number format/style.
function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

$ws.Range('D2').Value = '55.255.56'
$ws.Range('E2').Value = '  -4.73%  '
$ws.Range('D3').Value = '2.932.67'
$ws.Range('E3').Value = '  -7.57%  '
Set-TextValue $ws.Range('D4') "0.999"
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue $ws.Range('D5') "478.47"
$ws.Range('E5').Value = '  -9.50%  '
Set-TextValue $ws.Range('D6') "129.13"
$ws.Range('E6').Value = '  -3.94%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '2.934.38'
$ws.Range('E8').Value = '  -7.47%  '
Set-TextValue $ws.Range('D9') "0.414"
$ws.Range('E9').Value = '  -8.66%  '
Set-TextValue $ws.Range('D10') "6.87"
$ws.Range('E10').Value = '  -6.66%  '
Set-TextValue $ws.Range('D11') "0.0994"
$ws.Range('E11').Value = '  -11.03%  '
Set-TextValue $ws.Range('D12') "0.344"
$ws.Range('E12').Value = '  -12.39%  '
Set-TextValue $ws.Range('D13') "0.125"
$ws.Range('E13').Value = '  -1.59%  '
$ws.Range('D14').Value = '3.425.24'
$ws.Range('E14').Value = '  -7.77%  '
Set-TextValue $ws.Range('D15') "24.08"
$ws.Range('E15').Value = '  -7.17%  '
$ws.Range('D16').Value = '55.149.93'
$ws.Range('E16').Value = '  -4.82%  '
$ws.Range('D17').Value = '2.926.64'
$ws.Range('E17').Value = '  -7.49%  '
Set-TextValue $ws.Range('D18') "0.0000137"
$ws.Range('E18').Value = '  -10.88%  '
Set-TextValue $ws.Range('D19') "5.51"
$ws.Range('E19').Value = '  -5.61%  '
Set-TextValue $ws.Range('D20') "11.79"
$ws.Range('E20').Value = '  -10.27%  '
Set-TextValue $ws.Range('D21') "7.32"
$ws.Range('E21').Value = '  -9.80%  '
Set-TextValue $ws.Range('D22') "308.15"
$ws.Range('E22').Value = '  -11.72%  '
$ws.Range('E23').Value = '  -0.22%  '
Set-TextValue $ws.Range('D24') "0.454"
$ws.Range('E24').Value = '  -11.61%  '
Set-TextValue $ws.Range('D25') "59.34"
$ws.Range('E25').Value = '  -15.00%  '
$ws.Range('E26').Value = '  +0.08%  '
Set-TextValue $ws.Range('D27') "0.155"
$ws.Range('E27').Value = '  -7.32%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '0.0₃0833'
$ws.Range('E29').Value = '  -13.93%  '
Set-TextValue $ws.Range('D30') "6.40"
$ws.Range('E30').Value = '  -7.62%  '
Set-TextValue $ws.Range('D31') "6.41"
$ws.Range('E31').Value = '  -8.06%  '
Set-TextValue $ws.Range('D32') "1.15"
$ws.Range('E32').Value = '  -6.23%  '
Set-TextValue $ws.Range('D33') "1.66"
$ws.Range('E33').Value = '  -12.28%  '
Set-TextValue $ws.Range('D34') "19.07"
$ws.Range('E34').Value = '  -12.70%  '
Set-TextValue $ws.Range('D35') "147.14"
$ws.Range('E35').Value = '  -8.41%  '
Set-TextValue $ws.Range('D36') "4.28"
$ws.Range('E36').Value = '  -13.28%  '
Set-TextValue $ws.Range('D37') "5.54"
$ws.Range('E37').Value = '  -11.94%  '
Set-TextValue $ws.Range('D38') "1.25"
$ws.Range('E38').Value = '  -11.72%  '
Set-TextValue $ws.Range('D39') "23.37"
$ws.Range('E39').Value = '  -10.12%  '
Set-TextValue $ws.Range('D40') "0.0634"
$ws.Range('E40').Value = '  -9.43%  '
$ws.Range('D41').Value = '2.958.85'
$ws.Range('E41').Value = '  -7.53%  '
Set-TextValue $ws.Range('D42') "0.999"
$ws.Range('E42').Value = '  -0.17%  '
Set-TextValue $ws.Range('D43') "35.84"
$ws.Range('E43').Value = '  -12.03%  '
Set-TextValue $ws.Range('D44') "0.987"
$ws.Range('E44').Value = '  -9.75%  '
Set-TextValue $ws.Range('D45') "0.622"
$ws.Range('E45').Value = '  -11.42%  '
Set-TextValue $ws.Range('D46') "1.36"
$ws.Range('E46').Value = '  -7.75%  '
Set-TextValue $ws.Range('D47') "3.47"
$ws.Range('E47').Value = '  -12.57%  '
$ws.Range('D48').Value = '2.079.17'
$ws.Range('E48').Value = '  -8.62%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D49') "0.0225"
$ws.Range('E49').Value = '  -5.12%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D50') "5.50"
$ws.Range('E50').Value = '  -11.50%  '
Set-TextValue $ws.Range('D51') "18.51"
$ws.Range('E51').Value = '  -10.56%  '
